# Update database: shift the rolling 12-month columns forward one year
# (drop the oldest 1396/12 column's data, add a new 1401/12 column) and
# refresh the underlying figures - per "update database and change
# read_price algorithm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column headers (row 8 and row 24) -------------------------------
# Old: E=1396/12 F=1397/12 G=1398/12 H=1399/12 I=1400/12
# New: E=1397/12 F=1398/12 G=1399/12 H=1400/12 I=1401/12
$ws.Range("E8").Value  = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value  = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value  = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value  = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value  = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E24").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# --- Row 17: هزینه حقوق و دستمزد --------------------------------------
$ws.Range("E17").Value = 2595596
$ws.Range("F17").Value = 662746
$ws.Range("G17").Value = 1076749
$ws.Range("H17").Value = 1795832
$ws.Range("I17").Value = 3618828

# --- Row 19: سایر هزینه ها ---------------------------------------------
$ws.Range("E19").Value = 1713110
$ws.Range("F19").Value = 3048205
$ws.Range("G19").Value = 8823585
$ws.Range("H19").Value = 8471786
$ws.Range("I19").Value = 10394897

# --- Row 20: جمع ---------------------------------------------------------
$ws.Range("E20").Value = 4308706
$ws.Range("F20").Value = 3710951
$ws.Range("G20").Value = 9900334
$ws.Range("H20").Value = 10267618
$ws.Range("I20").Value = 14013725

# --- Row 26: تعداد پرسنل غیر تولیدی شرکت --------------------------------
$ws.Range("E26").Value = 821
$ws.Range("F26").Value = 789
$ws.Range("G26").Value = 790
$ws.Range("H26").Value = 756
$ws.Range("I26").Value = 739

# --- Row 27: تعداد پرسنل تولیدی شرکت -------------------------------------
$ws.Range("E27").Value = 1667
$ws.Range("F27").Value = 1613
$ws.Range("G27").Value = 1556
$ws.Range("H27").Value = 1520
$ws.Range("I27").Value = 1487
